$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.286832544864788;  C = 10.34677158129881;  D = 0.1494219747398047;  E = 10.19245300693656;  G = 23.97547910783996 }
    3  = @{ B = 1.455362044514542;  C = 1.655778082260271;   D = 0.1494219747398047;  E = 0.4942365360607697; G = 3.754798637575387 }
    4  = @{ B = 1.455362044514542;  C = 0.306821227259698;   D = 0.7527432677738641;  E = 0.4942365360607697; G = 3.009163075608874 }
    5  = @{ B = 1.455362044514542;  C = 1.655778082260271;   D = 0.7527432677738641;  E = 0.4942365360607697; G = 4.358119930609447 }
    6  = @{ B = 3.286832544864788;  C = 1.655778082260271;   D = 0.7527432677738641;  E = 0.4942365360607697; G = 6.189590430959694 }
    7  = @{ B = 0.2917716402565462; C = 0.306821227259698;   D = 0.7527432677738641;  E = 0.4942365360607697; G = 1.845572671350878 }
    8  = @{ B = 0.6606524410359556; C = 0.306821227259698;   D = 0.7527432677738641;  E = 1133.036916526867;  G = 1134.757133462937 }
    9  = @{ B = 0.1190320826869504; C = 10.34677158129881;   D = 261.3203778131603;   E = 10.19245300693656;  G = 281.9786344840826 }
    10 = @{ B = 1.455362044514542;  C = 1.655778082260271;   D = 3.537761648806719;   E = 10.19245300693656;  G = 16.84135478251809 }
    11 = @{ B = 3.286832544864788;  C = 10.34677158129881;   D = 3.537761648806719;   E = 1133.036916526867;  G = 1150.208282301838 }
    12 = @{ B = 1.455362044514542;  C = 1.655778082260271;   D = 0.1494219747398047;  E = 0.4942365360607697; G = 3.754798637575387 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}

$wb.Save()
